$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width fixes ---------------------------------------------------
# col B: 8.7109375 -> 7.7109375 (one "character" narrower)
# col D: 13.7109375 -> 11.7109375 (two "characters" narrower)
$ws.Columns.Item(2).ColumnWidth = 6.83
$ws.Columns.Item(4).ColumnWidth = 10.83

# --- Data fix --------------------------------------------------------------
# Every data row except row 13 had been written with the wrong per-country
# file/currency/sum values (the "bibliotheca sql write error"). Row 13 held
# the correct values ("87811004_1121_GB" / 400 / "GBP" / 1309.43 / "1309.43"),
# so re-apply that row onto every other data row. Copy+PasteSpecial (rather
# than re-typing the values) keeps each destination cell's exact type, e.g.
# column E stays text, matching row 13 exactly.
$src = $ws.Range("A13:E13")
$src.Copy()
for ($r = 2; $r -le 23; $r++) {
    if ($r -eq 13) { continue }
    $dst = $ws.Range("A" + $r + ":E" + $r)
    $dst.PasteSpecial()
}
$excel.CutCopyMode = $false
